# Update weekly fruit/vegetable price data (Frambuesa - Vega Monumental Concepción)
# Rows 2-3 take the date/price values that previously belonged to rows 4-5,
# and rows 4-5 take the date/price values that previously belonged to rows 2-3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was 2021-02-08 / 44216 -> now 2021-01-18 / 44195)
$ws.Range("D2").Value = 44195
$ws.Range("N2").Value = 3000
$ws.Range("O2").Value = 3500
$ws.Range("P2").Value = 3250
$ws.Range("S2").Value = 1625

# Row 3 (was 2021-02-08 / 44216 -> now 2021-01-18 / 44195)
$ws.Range("D3").Value = 44195
$ws.Range("N3").Value = 2500
$ws.Range("O3").Value = 2500
$ws.Range("P3").Value = 2500
$ws.Range("S3").Value = 1250

# Row 4 (was 2021-01-18 / 44195 -> now 2021-02-08 / 44216)
$ws.Range("D4").Value = 44216
$ws.Range("N4").Value = 3500
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3750
$ws.Range("S4").Value = 1875

# Row 5 (was 2021-01-18 / 44195 -> now 2021-02-08 / 44216)
$ws.Range("D5").Value = 44216
$ws.Range("N5").Value = 3000
$ws.Range("O5").Value = 3000
$ws.Range("P5").Value = 3000
$ws.Range("S5").Value = 1500
